$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.721.34"
Set-TextValue $ws.Range("E2") "  -4.25%  "
Set-TextValue $ws.Range("D3") "1.817.00"
Set-TextValue $ws.Range("E3") "  -3.11%  "
Set-TextValue $ws.Range("E4") "  -0.09%  "
Set-TextValue $ws.Range("D5") "277.29"
Set-TextValue $ws.Range("E5") "  -7.97%  "
Set-TextValue $ws.Range("E6") "  -0.09%  "
Set-TextValue $ws.Range("E7") "  -5.35%  "
Set-TextValue $ws.Range("D8") "0.3528"
Set-TextValue $ws.Range("E8") "  -6.13%  "
Set-TextValue $ws.Range("D9") "44.33"
Set-TextValue $ws.Range("E9") "  -2.39%  "
Set-TextValue $ws.Range("D10") "0.06660"
Set-TextValue $ws.Range("E10") "  -7.34%  "
Set-TextValue $ws.Range("D11") "20.08"
Set-TextValue $ws.Range("E11") "  -7.08%  "
Set-TextValue $ws.Range("D12") "0.8262"
Set-TextValue $ws.Range("E12") "  -7.14%  "
Set-TextValue $ws.Range("D13") "0.07875"
Set-TextValue $ws.Range("E13") "  -3.39%  "
Set-TextValue $ws.Range("D14") "1.812.27"
Set-TextValue $ws.Range("E14") "  -3.46%  "
Set-TextValue $ws.Range("D15") "5.070"
Set-TextValue $ws.Range("E15") "  -4.62%  "
Set-TextValue $ws.Range("D16") "87.58"
Set-TextValue $ws.Range("E16") "  -6.39%  "
Set-TextValue $ws.Range("D17") "1.000"
Set-TextValue $ws.Range("E17") "  -0.17%  "
Set-TextValue $ws.Range("E18") "  -5.06%  "
Set-TextValue $ws.Range("E19") "  -0.04%  "
Set-TextValue $ws.Range("D20") "0.000008017"
Set-TextValue $ws.Range("E20") "  -6.24%  "
Set-TextValue $ws.Range("D21") "25.761.64"
Set-TextValue $ws.Range("E21") "  -4.23%  "
Set-TextValue $ws.Range("D22") "4.740"
Set-TextValue $ws.Range("E22") "  -4.97%  "
Set-TextValue $ws.Range("D23") "9.999"
Set-TextValue $ws.Range("E23") "  -6.04%  "
Set-TextValue $ws.Range("E24") "  -4.94%  "
Set-TextValue $ws.Range("D25") "142.28"
Set-TextValue $ws.Range("E25") "  -2.80%  "
Set-TextValue $ws.Range("D26") "2.207"
Set-TextValue $ws.Range("E26") "  -3.92%  "
Set-TextValue $ws.Range("D27") "1.676"
Set-TextValue $ws.Range("E27") "  -3.16%  "
Set-TextValue $ws.Range("D28") "17.11"
Set-TextValue $ws.Range("E28") "  -5.44%  "
Set-TextValue $ws.Range("D29") "109.71"
Set-TextValue $ws.Range("E29") "  -3.81%  "
Set-TextValue $ws.Range("D30") "4.339"
Set-TextValue $ws.Range("E30") "  -8.19%  "
Set-TextValue $ws.Range("D31") "4.233"
Set-TextValue $ws.Range("E31") "  -8.39%  "
Set-TextValue $ws.Range("D32") "0.08778"
Set-TextValue $ws.Range("E32") "  -4.06%  "
Set-TextValue $ws.Range("D33") "0.04882"
Set-TextValue $ws.Range("E33") "  -2.50%  "
Set-TextValue $ws.Range("D34") "0.7271"
Set-TextValue $ws.Range("E34") "  -10.76%  "
Set-TextValue $ws.Range("D35") "1.139"
Set-TextValue $ws.Range("E35") "  -3.25%  "
Set-TextValue $ws.Range("D36") "2.869"
Set-TextValue $ws.Range("E36") "  -2.58%  "
Set-TextValue $ws.Range("B37") "Frax"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D37") "0.9997"
Set-TextValue $ws.Range("E37") "  -0.17%  "
Set-TextValue $ws.Range("B38") "MXToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "3.126"
Set-TextValue $ws.Range("E38") "  -2.67%  "
Set-TextValue $ws.Range("B39") "RenderToken"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.371"
Set-TextValue $ws.Range("E39") "  -9.35%  "
Set-TextValue $ws.Range("B40") "VeChain"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.01851"
Set-TextValue $ws.Range("E40") "  -5.17%  "
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5170"
Set-TextValue $ws.Range("E41") "  -14.60%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.9638"
Set-TextValue $ws.Range("E42") "  -10.02%  "
Set-TextValue $ws.Range("B43") "FraxShare"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "6.218"
Set-TextValue $ws.Range("E43") "  -6.31%  "
Set-TextValue $ws.Range("B44") "Quant"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D44") "110.42"
Set-TextValue $ws.Range("E44") "  -4.16%  "
Set-TextValue $ws.Range("B45") "Aptos"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D45") "8.025"
Set-TextValue $ws.Range("E45") "  -10.11%  "
Set-TextValue $ws.Range("B46") "PaxDollar"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "1.001"
Set-TextValue $ws.Range("E46") "  -0.07%  "
Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.4553"
Set-TextValue $ws.Range("E47") "  -10.89%  "
Set-TextValue $ws.Range("B48") "Algorand"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D48") "0.1364"
Set-TextValue $ws.Range("E48") "  -8.61%  "
Set-TextValue $ws.Range("B49") "Elrond"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D49") "36.55"
Set-TextValue $ws.Range("E49") "  -3.05%  "
Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "9.241"
Set-TextValue $ws.Range("E50") "  -7.14%  "
Set-TextValue $ws.Range("B51") "NEARProtocol"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.497"
Set-TextValue $ws.Range("E51") "  -8.29%  "
